$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") '28.345.69'
$ws.Range("E2").Value = '  +5.80%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.810.49'
$ws.Range("E3").Value = '  +5.22%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.002'
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
Set-TextValue $ws.Range("D5") '318.20'
$ws.Range("E5").Value = '  +3.07%  '

# Row 6
Set-TextValue $ws.Range("D6") '1.002'
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.5708'
$ws.Range("E7").Value = '  +17.08%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3884'
$ws.Range("E8").Value = '  +11.31%  '

# Row 9
Set-TextValue $ws.Range("D9") '43.12'
$ws.Range("E9").Value = '  +1.03%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.07593'
$ws.Range("E10").Value = '  +5.22%  '

# Row 11
Set-TextValue $ws.Range("D11") '1.137'
$ws.Range("E11").Value = '  +8.66%  '

# Row 12
Set-TextValue $ws.Range("D12") '1.002'
$ws.Range("E12").Value = '  +0.03%  '

# Row 13
Set-TextValue $ws.Range("D13") '21.17'
$ws.Range("E13").Value = '  +7.12%  '

# Row 14
Set-TextValue $ws.Range("D14") '6.253'
$ws.Range("E14").Value = '  +6.87%  '

# Row 15
Set-TextValue $ws.Range("D15") '1.809.98'
$ws.Range("E15").Value = '  +5.43%  '

# Row 16
Set-TextValue $ws.Range("D16") '7.249'
$ws.Range("E16").Value = '  +6.81%  '

# Row 17
Set-TextValue $ws.Range("D17") '91.96'
$ws.Range("E17").Value = '  +6.65%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.00001073'
$ws.Range("E18").Value = '  +3.93%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.06475'
$ws.Range("E19").Value = '  +1.39%  '

# Row 20
Set-TextValue $ws.Range("D20") '1.002'
$ws.Range("E20").Value = '  +0.04%  '

# Row 21
$ws.Range("E21").Value = '  +5.12%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.999'
$ws.Range("E22").Value = '  +5.21%  '

# Row 23
Set-TextValue $ws.Range("D23") '28.358.76'
$ws.Range("E23").Value = '  +5.55%  '

# Row 24
Set-TextValue $ws.Range("D24") '11.28'
$ws.Range("E24").Value = '  +3.36%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.135'
$ws.Range("E25").Value = '  +4.33%  '

# Row 26
Set-TextValue $ws.Range("D26") '158.16'
$ws.Range("E26").Value = '  +2.53%  '

# Row 27
Set-TextValue $ws.Range("D27") '20.80'
$ws.Range("E27").Value = '  +5.19%  '

# Row 28
Set-TextValue $ws.Range("D28") '2.439'
$ws.Range("E28").Value = '  +18.53%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.017.56'
$ws.Range("E29").Value = '  +5.37%  '

# Row 30
Set-TextValue $ws.Range("D30") '123.75'
$ws.Range("E30").Value = '  +3.50%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.157'
$ws.Range("E31").Value = '  +11.62%  '

# Row 32
$ws.Range("E32").Value = '  +13.62%  '

# Row 33
Set-TextValue $ws.Range("D33") '5.774'
$ws.Range("E33").Value = '  +8.32%  '

# Row 34
Set-TextValue $ws.Range("D34") '3.639'
$ws.Range("E34").Value = '  +1.75%  '

# Row 35
$ws.Range("B35").Value = 'Algorand'
$ws.Range("C35").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D35") '0.2204'
$ws.Range("E35").Value = '  +11.43%  '

# Row 36
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D36") '8.919'
$ws.Range("E36").Value = '  +20.47%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D37") '0.02321'
$ws.Range("E37").Value = '  +7.24%  '

# Row 38
Set-TextValue $ws.Range("D38") '11.66'
$ws.Range("E38").Value = '  +7.03%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.06117'
$ws.Range("E39").Value = '  +4.25%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.6402'
$ws.Range("E40").Value = '  +7.61%  '

# Row 41
Set-TextValue $ws.Range("D41") '5.036'
$ws.Range("E41").Value = '  +6.38%  '

# Row 42
$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D42") '1.000'
$ws.Range("E42").Value = '  -0.01%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D43") '1.161'
$ws.Range("E43").Value = '  +4.24%  '

# Row 44
Set-TextValue $ws.Range("D44") '1.381'
$ws.Range("E44").Value = '  -2.97%  '

# Row 45
Set-TextValue $ws.Range("D45") '13.36'
$ws.Range("E45").Value = '  +4.85%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.6000'
$ws.Range("E46").Value = '  +7.46%  '

# Row 47
Set-TextValue $ws.Range("D47") '3.702'
$ws.Range("E47").Value = '  +3.72%  '

# Row 48
Set-TextValue $ws.Range("D48") '121.95'
$ws.Range("E48").Value = '  +2.28%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.941'
$ws.Range("E49").Value = '  +6.20%  '

# Row 50
$ws.Range("E50").Value = '  +5.45%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.06861'
$ws.Range("E51").Value = '  +3.56%  '
